$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = 22047
$ws.Range("E19").Value = 32930047
$ws.Range("C37").Value = 33662
$ws.Range("D37").Value = 10428
$ws.Range("E37").Value = 49422067
$ws.Range("C55").Value = 58152
$ws.Range("E55").Value = 91763268
$ws.Range("C66").Value = 70589
$ws.Range("E66").Value = 105647624
$ws.Range("E78").Value = 57752877
$ws.Range("C86").Value = 33635
$ws.Range("D86").Value = 10339
$ws.Range("E86").Value = 50614523
$ws.Range("E96").Value = 195863097
$ws.Range("C109").Value = 47406
$ws.Range("E109").Value = 78515227
$ws.Range("C119").Value = 23970
$ws.Range("E119").Value = 35374148
$ws.Range("C152").Value = 8674
$ws.Range("E152").Value = 13787642
$ws.Range("C159").Value = 25079
$ws.Range("E159").Value = 39638350
$ws.Range("C198").Value = 43323
$ws.Range("E198").Value = 66251883
$ws.Range("C218").Value = 11756
$ws.Range("E218").Value = 18688473
$ws.Range("C278").Value = 16971
$ws.Range("E278").Value = 27498705
$ws.Range("C296").Value = 12511
$ws.Range("E296").Value = 22202911
$ws.Range("D305").Value = 15228
$ws.Range("C324").Value = 28698
$ws.Range("E324").Value = 42889341
$ws.Range("C360").Value = 33114
$ws.Range("E360").Value = 49889741
$ws.Range("C401").Value = 97757
$ws.Range("E401").Value = 145934805
$ws.Range("C412").Value = 29328
$ws.Range("E412").Value = 45383476
$ws.Range("C438").Value = 59716
$ws.Range("D438").Value = 15779
$ws.Range("E438").Value = 96956189
$ws.Range("C446").Value = 108992
$ws.Range("D446").Value = 26324
$ws.Range("E446").Value = 180641979
$ws.Range("C455").Value = 265634
$ws.Range("D455").Value = 67475
$ws.Range("E455").Value = 421218549
$ws.Range("C472").Value = 68580
$ws.Range("D472").Value = 18408
$ws.Range("E472").Value = 108736485
$ws.Range("C481").Value = 121760
$ws.Range("D481").Value = 26742
$ws.Range("E481").Value = 212448611
$ws.Range("C490").Value = 94214
$ws.Range("E490").Value = 153842183
$ws.Range("C499").Value = 71741
$ws.Range("E499").Value = 120296901
$ws.Range("C507").Value = 66755
$ws.Range("D507").Value = 18643
$ws.Range("E507").Value = 105583800
$ws.Range("C540").Value = 33560
$ws.Range("D540").Value = 10534
$ws.Range("E540").Value = 51736035
$ws.Range("C565").Value = 11319
$ws.Range("D565").Value = 3573
$ws.Range("E565").Value = 18206075
$ws.Range("C572").Value = 44267
$ws.Range("D572").Value = 13805
$ws.Range("E572").Value = 67917933
$ws.Range("C590").Value = 41002
$ws.Range("E590").Value = 60347867
$ws.Range("C631").Value = 111160
$ws.Range("E631").Value = 166159428
$ws.Range("C652").Value = 26090
$ws.Range("E652").Value = 38444729
$ws.Range("C669").Value = 44259
$ws.Range("D669").Value = 14320
$ws.Range("E669").Value = 64155643
$ws.Range("C687").Value = 10284
$ws.Range("D687").Value = 3368
$ws.Range("E687").Value = 14072791
$ws.Range("C693").Value = 25101
$ws.Range("D693").Value = 7662
$ws.Range("E693").Value = 36647614
$ws.Range("C710").Value = 50172
$ws.Range("E710").Value = 71962122
$ws.Range("C726").Value = 91327
$ws.Range("E726").Value = 133670667
$ws.Range("E749").Value = 126256866
$ws.Range("C759").Value = 11625
$ws.Range("D759").Value = 3763
$ws.Range("E759").Value = 17284624
$ws.Range("C765").Value = 5407
$ws.Range("E765").Value = 8778745
$ws.Range("C797").Value = 66460
$ws.Range("D797").Value = 21675
$ws.Range("E797").Value = 100767215
$ws.Range("C820").Value = 19697
$ws.Range("D820").Value = 6165
$ws.Range("E820").Value = 30618101
$ws.Range("C838").Value = 13109
$ws.Range("D838").Value = 4046
$ws.Range("E838").Value = 20176940
$ws.Range("C847").Value = 108076
$ws.Range("E847").Value = 164182835
$ws.Range("C861").Value = 134231
$ws.Range("E861").Value = 201915733
$ws.Range("C874").Value = 18236
$ws.Range("E874").Value = 29210227
$ws.Range("C885").Value = 85283
$ws.Range("D885").Value = 26091
$ws.Range("E885").Value = 125689469
$ws.Range("C895").Value = 43340
$ws.Range("E895").Value = 63854196
